$wb = $excel.ActiveWorkbook

# ---- Sheet ALC: 96 cell updates ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1327.7894
$ws.Range("J19").Value = 1498.9286
$ws.Range("L19").Value = 1498.9286
$ws.Range("N19").Value = -1848.9286
$ws.Range("H28").Value = 3957.923
$ws.Range("I28").Value = 1057.6666
$ws.Range("J28").Value = 6443.857
$ws.Range("K28").Value = 1057.6666
$ws.Range("L28").Value = 6443.857
$ws.Range("M28").Value = -572.6666
$ws.Range("N28").Value = -7413.857
$ws.Range("H40").Value = 4123.8774
$ws.Range("I40").Value = 2448.8
$ws.Range("K40").Value = 2448.8
$ws.Range("M40").Value = -2273.8
$ws.Range("H64").Value = 12343.823
$ws.Range("J64").Value = 6872.6924
$ws.Range("L64").Value = 6872.6924
$ws.Range("N64").Value = -7368.6924
$ws.Range("H67").Value = 12343.823
$ws.Range("J67").Value = 6872.6924
$ws.Range("L67").Value = 6872.6924
$ws.Range("N67").Value = -8588.6924
$ws.Range("H74").Value = 7244.16
$ws.Range("I74").Value = 4280.8
$ws.Range("K74").Value = 4280.8
$ws.Range("M74").Value = -3344.8
$ws.Range("H77").Value = 7244.16
$ws.Range("I77").Value = 4280.8
$ws.Range("K77").Value = 21404
$ws.Range("M77").Value = -16724
$ws.Range("H96").Value = 322850.66
$ws.Range("I96").Value = 825.2
$ws.Range("J96").Value = 725382.5
$ws.Range("K96").Value = 2475.6
$ws.Range("L96").Value = 2176147.5
$ws.Range("M96").Value = -1102.6
$ws.Range("N96").Value = -2178893.5
$ws.Range("H100").Value = 963.9
$ws.Range("I100").Value = 698.8421
$ws.Range("K100").Value = 698.8421
$ws.Range("M100").Value = -157.8421
$ws.Range("H107").Value = 18574558
$ws.Range("I107").Value = 19667090
$ws.Range("K107").Value = 19667090
$ws.Range("M107").Value = -19665170
$ws.Range("H112").Value = 6504.087
$ws.Range("I112").Value = 2114
$ws.Range("K112").Value = 6342
$ws.Range("M112").Value = -5234
$ws.Range("H116").Value = 4983.737
$ws.Range("I116").Value = 4421.6665
$ws.Range("J116").Value = 5243.154
$ws.Range("K116").Value = 4421.6665
$ws.Range("L116").Value = 5243.154
$ws.Range("M116").Value = -979.6665000000003
$ws.Range("N116").Value = -12127.154
$ws.Range("H127").Value = 2400
$ws.Range("I127").Value = 2400
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 7200
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = -2240
$ws.Range("N127").ClearContents()
$ws.Range("H132").Value = 21280044
$ws.Range("I132").Value = 26318914
$ws.Range("J132").Value = 4820.6665
$ws.Range("K132").Value = 78956742
$ws.Range("L132").Value = 14461.9995
$ws.Range("M132").Value = -78954212
$ws.Range("N132").Value = -19521.9995
$ws.Range("H135").Value = 1098.1154
$ws.Range("I135").Value = 631.875
$ws.Range("J135").Value = 1844.1
$ws.Range("K135").Value = 5686.875
$ws.Range("L135").Value = 16596.9
$ws.Range("M135").Value = -3151.875
$ws.Range("N135").Value = -21666.9
$ws.Range("H137").Value = 75398.625
$ws.Range("I137").Value = 75398.625
$ws.Range("K137").Value = 226195.875
$ws.Range("M137").Value = -223645.875
$ws.Range("H138").Value = 4034.3618
$ws.Range("I138").Value = 3110.8076
$ws.Range("J138").Value = 5177.8096
$ws.Range("K138").Value = 9332.4228
$ws.Range("L138").Value = 15533.4288
$ws.Range("M138").Value = -4192.4228
$ws.Range("N138").Value = -25813.4288
$ws.Range("H141").Value = 9323.459000000001
$ws.Range("I141").Value = 3895.9062
$ws.Range("J141").Value = 44059.8
$ws.Range("K141").Value = 11687.7186
$ws.Range("L141").Value = 132179.4
$ws.Range("M141").Value = -6507.7186
$ws.Range("N141").Value = -142539.4

# ---- Sheet ARM: 71 cell updates ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2572336.8
$ws.Range("I2").Value = 3328200.8
$ws.Range("J2").Value = 2398.8
$ws.Range("K2").Value = 3328200.8
$ws.Range("L2").Value = 2398.8
$ws.Range("M2").Value = -3328087.8
$ws.Range("N2").Value = -2624.8
$ws.Range("H32").Value = 3393.7231
$ws.Range("I32").Value = 2047.3455
$ws.Range("K32").Value = 2047.3455
$ws.Range("M32").Value = -1760.3455
$ws.Range("H45").Value = 4498986
$ws.Range("I45").Value = 6851441
$ws.Range("J45").Value = 7935.4546
$ws.Range("K45").Value = 6851441
$ws.Range("L45").Value = 7935.4546
$ws.Range("M45").Value = -6851064
$ws.Range("N45").Value = -8689.454600000001
$ws.Range("H61").Value = 5740.6055
$ws.Range("I61").Value = 5760.757
$ws.Range("J61").Value = 4995
$ws.Range("K61").Value = 5760.757
$ws.Range("L61").Value = 4995
$ws.Range("M61").Value = -5548.757
$ws.Range("N61").Value = -5419
$ws.Range("H74").Value = 36744.934
$ws.Range("I74").Value = 7364.1875
$ws.Range("J74").Value = 154267.92
$ws.Range("K74").Value = 7364.1875
$ws.Range("L74").Value = 154267.92
$ws.Range("M74").Value = -6490.1875
$ws.Range("N74").Value = -156015.92
$ws.Range("H77").Value = 36744.934
$ws.Range("I77").Value = 7364.1875
$ws.Range("J77").Value = 154267.92
$ws.Range("K77").Value = 36820.9375
$ws.Range("L77").Value = 771339.6000000001
$ws.Range("M77").Value = -32452.9375
$ws.Range("N77").Value = -780075.6000000001
$ws.Range("H97").Value = 1546120
$ws.Range("I97").Value = 2157718.5
$ws.Range("J97").Value = 17124
$ws.Range("K97").Value = 2157718.5
$ws.Range("L97").Value = 17124
$ws.Range("M97").Value = -2157222.5
$ws.Range("N97").Value = -18116
$ws.Range("H116").Value = 2572336.8
$ws.Range("I116").Value = 3328200.8
$ws.Range("J116").Value = 2398.8
$ws.Range("K116").Value = 3328200.8
$ws.Range("L116").Value = 2398.8
$ws.Range("M116").Value = -3325906.8
$ws.Range("N116").Value = -6986.8
$ws.Range("H122").Value = 615904
$ws.Range("I122").Value = 3013.6
$ws.Range("J122").Value = 1491461.8
$ws.Range("K122").Value = 9040.799999999999
$ws.Range("L122").Value = 4474385.4
$ws.Range("M122").Value = -6590.799999999999
$ws.Range("N122").Value = -4479285.4
$ws.Range("H132").Value = 5124.61
$ws.Range("I132").Value = 4634.4062
$ws.Range("K132").Value = 13903.2186
$ws.Range("M132").Value = -11373.2186
$ws.Range("H136").Value = 5740.6055
$ws.Range("I136").Value = 5760.757
$ws.Range("J136").Value = 4995
$ws.Range("K136").Value = 17282.271
$ws.Range("L136").Value = 14985
$ws.Range("M136").Value = -14732.271
$ws.Range("N136").Value = -20085

# ---- Sheet BSM: 68 cell updates ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2572336.8
$ws.Range("I3").Value = 3328200.8
$ws.Range("J3").Value = 2398.8
$ws.Range("K3").Value = 3328200.8
$ws.Range("L3").Value = 2398.8
$ws.Range("M3").Value = -3328086.8
$ws.Range("N3").Value = -2626.8
$ws.Range("H15").Value = 8999.5
$ws.Range("J15").Value = 8999.5
$ws.Range("L15").Value = 8999.5
$ws.Range("N15").Value = -9453.5
$ws.Range("H35").Value = 31000
$ws.Range("J35").Value = 31000
$ws.Range("L35").Value = 31000
$ws.Range("N35").Value = -31620
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H86").Value = 3855801.5
$ws.Range("I86").Value = 5567683
$ws.Range("J86").Value = 4067.75
$ws.Range("K86").Value = 5567683
$ws.Range("L86").Value = 4067.75
$ws.Range("M86").Value = -5566560
$ws.Range("N86").Value = -6313.75
$ws.Range("H89").Value = 3855801.5
$ws.Range("I89").Value = 5567683
$ws.Range("J89").Value = 4067.75
$ws.Range("K89").Value = 27838415
$ws.Range("L89").Value = 20338.75
$ws.Range("M89").Value = -27832799
$ws.Range("N89").Value = -31570.75
$ws.Range("H105").Value = 3343140
$ws.Range("I105").Value = 4234325
$ws.Range("K105").Value = 4234325
$ws.Range("M105").Value = -4232578
$ws.Range("H107").Value = 2465160.5
$ws.Range("I107").Value = 3248299.5
$ws.Range("J107").Value = 3866.1428
$ws.Range("K107").Value = 3248299.5
$ws.Range("L107").Value = 3866.1428
$ws.Range("M107").Value = -3246379.5
$ws.Range("N107").Value = -7706.1428
$ws.Range("H132").Value = 51421.43
$ws.Range("J132").Value = 51421.43
$ws.Range("L132").Value = 51421.43
$ws.Range("N132").Value = -61541.43
$ws.Range("H134").Value = 10421.706
$ws.Range("I134").Value = 12117.131
$ws.Range("K134").Value = 36351.393
$ws.Range("M134").Value = -33816.393
$ws.Range("H135").Value = 72694.75
$ws.Range("J135").Value = 72694.75
$ws.Range("L135").Value = 72694.75
$ws.Range("N135").Value = -82834.75
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ---- Sheet CRP: 104 cell updates ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 2002
$ws.Range("I11").Value = 1004.5
$ws.Range("J11").Value = 2999.5
$ws.Range("K11").Value = 1004.5
$ws.Range("L11").Value = 2999.5
$ws.Range("M11").Value = -864.5
$ws.Range("N11").Value = -3279.5
$ws.Range("H22").Value = 539
$ws.Range("I22").Value = 617.9091
$ws.Range("J22").Value = 249.66667
$ws.Range("K22").Value = 617.9091
$ws.Range("L22").Value = 249.66667
$ws.Range("M22").Value = -267.9091
$ws.Range("N22").Value = -949.6666700000001
$ws.Range("H31").Value = 15591.978
$ws.Range("I31").Value = 5391.609
$ws.Range("J31").Value = 19146.652
$ws.Range("K31").Value = 5391.609
$ws.Range("L31").Value = 19146.652
$ws.Range("M31").Value = -5096.609
$ws.Range("N31").Value = -19736.652
$ws.Range("H34").Value = 15591.978
$ws.Range("I34").Value = 5391.609
$ws.Range("J34").Value = 19146.652
$ws.Range("K34").Value = 5391.609
$ws.Range("L34").Value = 19146.652
$ws.Range("M34").Value = -5189.609
$ws.Range("N34").Value = -19550.652
$ws.Range("H35").Value = 1720
$ws.Range("I35").Value = 1720
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1720
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -1426
$ws.Range("N35").ClearContents()
$ws.Range("H50").Value = 6549.125
$ws.Range("J50").Value = 6549.125
$ws.Range("L50").Value = 6549.125
$ws.Range("N50").Value = -7799.125
$ws.Range("H58").Value = 2721.923
$ws.Range("I58").Value = 2048.75
$ws.Range("J58").Value = 3799
$ws.Range("K58").Value = 2048.75
$ws.Range("L58").Value = 3799
$ws.Range("M58").Value = -1845.75
$ws.Range("N58").Value = -4205
$ws.Range("H62").Value = 3386.4443
$ws.Range("I62").Value = 1496.3334
$ws.Range("J62").Value = 4331.5
$ws.Range("K62").Value = 1496.3334
$ws.Range("L62").Value = 4331.5
$ws.Range("M62").Value = -872.3334
$ws.Range("N62").Value = -5579.5
$ws.Range("H65").Value = 3386.4443
$ws.Range("I65").Value = 1496.3334
$ws.Range("J65").Value = 4331.5
$ws.Range("K65").Value = 7481.666999999999
$ws.Range("L65").Value = 21657.5
$ws.Range("M65").Value = -4361.666999999999
$ws.Range("N65").Value = -27897.5
$ws.Range("H105").Value = 1586.1578
$ws.Range("I105").Value = 1522.9231
$ws.Range("K105").Value = 1522.9231
$ws.Range("M105").Value = 224.0769
$ws.Range("H107").Value = 1393.95
$ws.Range("I107").Value = 1027.6875
$ws.Range("J107").Value = 2859
$ws.Range("K107").Value = 1027.6875
$ws.Range("L107").Value = 2859
$ws.Range("M107").Value = 892.3125
$ws.Range("N107").Value = -6699
$ws.Range("H122").Value = 3136.88
$ws.Range("I122").Value = 2039.6428
$ws.Range("J122").Value = 4533.364
$ws.Range("K122").Value = 6118.928400000001
$ws.Range("L122").Value = 13600.092
$ws.Range("M122").Value = -3668.928400000001
$ws.Range("N122").Value = -18500.092
$ws.Range("H132").Value = 51992.95
$ws.Range("I132").Value = 54672.05
$ws.Range("K132").Value = 164016.15
$ws.Range("M132").Value = -161486.15
$ws.Range("H134").Value = 7651.615
$ws.Range("I134").Value = 6504.4287
$ws.Range("J134").Value = 8990
$ws.Range("K134").Value = 19513.2861
$ws.Range("L134").Value = 26970
$ws.Range("M134").Value = -16978.2861
$ws.Range("N134").Value = -32040
$ws.Range("H136").Value = 2721.923
$ws.Range("I136").Value = 2048.75
$ws.Range("J136").Value = 3799
$ws.Range("K136").Value = 6146.25
$ws.Range("L136").Value = 11397
$ws.Range("M136").Value = -3596.25
$ws.Range("N136").Value = -16497
$ws.Range("H138").Value = 94966
$ws.Range("J138").Value = 94966
$ws.Range("L138").Value = 94966
$ws.Range("N138").Value = -105246
$ws.Range("H140").Value = 85680
$ws.Range("J140").Value = 85680
$ws.Range("L140").Value = 85680
$ws.Range("N140").Value = -96040

# ---- Sheet CUL: 161 cell updates ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 9976416
$ws.Range("J4").Value = 5244038
$ws.Range("L4").Value = 15732114
$ws.Range("N4").Value = -15732338
$ws.Range("H5").Value = 50350
$ws.Range("I5").Value = 695
$ws.Range("J5").Value = 100005
$ws.Range("K5").Value = 2085
$ws.Range("L5").Value = 300015
$ws.Range("M5").Value = -1973
$ws.Range("N5").Value = -300239
$ws.Range("H7").Value = 952.44446
$ws.Range("I7").Value = 1017.5789
$ws.Range("J7").Value = 797.75
$ws.Range("K7").Value = 3052.7367
$ws.Range("L7").Value = 2393.25
$ws.Range("M7").Value = -2940.7367
$ws.Range("N7").Value = -2617.25
$ws.Range("H10").Value = 53.4
$ws.Range("I10").Value = 53.4
$ws.Range("K10").Value = 160.2
$ws.Range("M10").Value = -21.19999999999999
$ws.Range("H12").Value = 152233.67
$ws.Range("I12").Value = 888888
$ws.Range("J12").Value = 4902.8
$ws.Range("K12").Value = 2666664
$ws.Range("L12").Value = 14708.4
$ws.Range("M12").Value = -2666491
$ws.Range("N12").Value = -15054.4
$ws.Range("H32").Value = 450000740
$ws.Range("I32").Value = 1000000000
$ws.Range("K32").Value = 3000000000
$ws.Range("M32").Value = -2999999717
$ws.Range("H34").Value = 1210.6666
$ws.Range("J34").Value = 2437.25
$ws.Range("L34").Value = 7311.75
$ws.Range("N34").Value = -7479.75
$ws.Range("H39").Value = 1649.8572
$ws.Range("I39").Value = 2874.5
$ws.Range("J39").Value = 1160
$ws.Range("K39").Value = 8623.5
$ws.Range("L39").Value = 3480
$ws.Range("M39").Value = -8329.5
$ws.Range("N39").Value = -4068
$ws.Range("H55").Value = 49396.953
$ws.Range("J55").Value = 171481.17
$ws.Range("L55").Value = 514443.51
$ws.Range("N55").Value = -514797.51
$ws.Range("H57").Value = 4794.3335
$ws.Range("J57").Value = 5049.875
$ws.Range("L57").Value = 15149.625
$ws.Range("N57").Value = -16267.625
$ws.Range("H60").Value = 2498.7273
$ws.Range("I60").Value = 2349.1
$ws.Range("K60").Value = 7047.299999999999
$ws.Range("M60").Value = -6796.299999999999
$ws.Range("H62").Value = 4497.5
$ws.Range("J62").Value = 4497.5
$ws.Range("L62").Value = 13492.5
$ws.Range("N62").Value = -14864.5
$ws.Range("H65").Value = 4497.5
$ws.Range("J65").Value = 4497.5
$ws.Range("L65").Value = 40477.5
$ws.Range("N65").Value = -47341.5
$ws.Range("H68").Value = 2062.0435
$ws.Range("I68").Value = 1832.5834
$ws.Range("K68").Value = 5497.7502
$ws.Range("M68").Value = -4686.7502
$ws.Range("H71").Value = 2062.0435
$ws.Range("I71").Value = 1832.5834
$ws.Range("K71").Value = 16493.2506
$ws.Range("M71").Value = -12437.2506
$ws.Range("H107").Value = 1131.0526
$ws.Range("I107").Value = 869.625
$ws.Range("J107").Value = 1321.1818
$ws.Range("K107").Value = 2608.875
$ws.Range("L107").Value = 3963.5454
$ws.Range("M107").Value = -688.875
$ws.Range("N107").Value = -7803.5454
$ws.Range("H110").Value = 20716.428
$ws.Range("I110").Value = 5000
$ws.Range("J110").Value = 21925.385
$ws.Range("K110").Value = 15000
$ws.Range("L110").Value = 65776.155
$ws.Range("M110").Value = -10910
$ws.Range("N110").Value = -73956.155
$ws.Range("H111").Value = 5109.6665
$ws.Range("I111").Value = 3449.5
$ws.Range("J111").Value = 8430
$ws.Range("K111").Value = 10348.5
$ws.Range("L111").Value = 25290
$ws.Range("M111").Value = -7281.5
$ws.Range("N111").Value = -31424
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").ClearContents()
$ws.Range("H114").Value = 1111533.5
$ws.Range("I114").Value = 453
$ws.Range("J114").Value = 5000315.5
$ws.Range("K114").Value = 1359
$ws.Range("L114").Value = 15000946.5
$ws.Range("M114").Value = 1895
$ws.Range("N114").Value = -15007454.5
$ws.Range("H115").Value = 4060.9167
$ws.Range("I115").Value = 3121.8333
$ws.Range("K115").Value = 9365.499899999999
$ws.Range("M115").Value = -8190.499899999999
$ws.Range("H118").Value = 2339.125
$ws.Range("I118").Value = 1953
$ws.Range("K118").Value = 5859
$ws.Range("M118").Value = -4616
$ws.Range("H119").Value = 5664.3335
$ws.Range("I119").Value = 996.5
$ws.Range("K119").Value = 2989.5
$ws.Range("M119").Value = 1848.5
$ws.Range("H120").Value = 18772.727
$ws.Range("J120").Value = 20150
$ws.Range("L120").Value = 60450
$ws.Range("N120").Value = -70126
$ws.Range("H123").Value = 10186.2
$ws.Range("I123").Value = 6999.6665
$ws.Range("K123").Value = 20998.9995
$ws.Range("M123").Value = -18548.9995
$ws.Range("H125").Value = 7099.3335
$ws.Range("I125").Value = 6132.5
$ws.Range("K125").Value = 18397.5
$ws.Range("M125").Value = -13477.5
$ws.Range("H126").Value = 30000
$ws.Range("J126").Value = 30000
$ws.Range("L126").Value = 90000
$ws.Range("N126").Value = -99880
$ws.Range("H129").Value = 910284.9
$ws.Range("I129").Value = 1250964
$ws.Range("K129").Value = 3752892
$ws.Range("M129").Value = -3747892
$ws.Range("H130").Value = 2541.625
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 2541.625
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 7624.875
$ws.Range("M130").ClearContents()
$ws.Range("N130").Value = -17664.875
$ws.Range("H131").Value = 3202.5
$ws.Range("I131").Value = 652.125
$ws.Range("K131").Value = 1956.375
$ws.Range("M131").Value = 3083.625
$ws.Range("H132").Value = 1049.75
$ws.Range("I132").Value = 1049.75
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9447.75
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6917.75
$ws.Range("N132").ClearContents()
$ws.Range("H135").Value = 50350
$ws.Range("I135").Value = 695
$ws.Range("J135").Value = 100005
$ws.Range("K135").Value = 6255
$ws.Range("L135").Value = 900045
$ws.Range("M135").Value = -3720
$ws.Range("N135").Value = -905115

# ---- Sheet GSM: 65 cell updates ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 10001
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H54").Value = 15000
$ws.Range("I54").Value = 15000
$ws.Range("K54").Value = 15000
$ws.Range("M54").Value = -14610
$ws.Range("H63").Value = 39499.5
$ws.Range("J63").Value = 39499.5
$ws.Range("L63").Value = 39499.5
$ws.Range("N63").Value = -40871.5
$ws.Range("H66").Value = 39499.5
$ws.Range("J66").Value = 39499.5
$ws.Range("L66").Value = 118498.5
$ws.Range("N66").Value = -125362.5
$ws.Range("H76").Value = 10001
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 10001
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H97").Value = 3408290.5
$ws.Range("I97").Value = 3975088.8
$ws.Range("K97").Value = 3975088.8
$ws.Range("M97").Value = -3974592.8
$ws.Range("H102").Value = 5666263
$ws.Range("I102").Value = 10103066
$ws.Range("K102").Value = 10103066
$ws.Range("M102").Value = -10101444
$ws.Range("H107").Value = 736.36365
$ws.Range("I107").Value = 833.3333
$ws.Range("K107").Value = 833.3333
$ws.Range("M107").Value = 1086.6667
$ws.Range("H113").Value = 9808009
$ws.Range("I113").Value = 20836396
$ws.Range("J113").Value = 4999.8887
$ws.Range("K113").Value = 20836396
$ws.Range("L113").Value = 4999.8887
$ws.Range("M113").Value = -20834226
$ws.Range("N113").Value = -9339.8887
$ws.Range("H122").Value = 1118849.2
$ws.Range("I122").Value = 1488140.6
$ws.Range("K122").Value = 4464421.800000001
$ws.Range("M122").Value = -4461971.800000001
$ws.Range("H126").Value = 5155070.5
$ws.Range("I126").Value = 2677042
$ws.Range("J126").Value = 10420881
$ws.Range("K126").Value = 8031126
$ws.Range("L126").Value = 31262643
$ws.Range("M126").Value = -8028656
$ws.Range("N126").Value = -31267583
$ws.Range("H132").Value = 7487.405
$ws.Range("I132").Value = 5731.4414
$ws.Range("J132").Value = 14950.25
$ws.Range("K132").Value = 17194.3242
$ws.Range("L132").Value = 44850.75
$ws.Range("M132").Value = -14664.3242
$ws.Range("N132").Value = -49910.75
$ws.Range("H140").Value = 61250
$ws.Range("J140").Value = 61250
$ws.Range("L140").Value = 61250
$ws.Range("N140").Value = -71610

# ---- Sheet LTW: 66 cell updates ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2123.9
$ws.Range("I16").Value = 1832.1111
$ws.Range("K16").Value = 1832.1111
$ws.Range("M16").Value = -1662.1111
$ws.Range("H25").Value = 26665
$ws.Range("I25").Value = 26665
$ws.Range("K25").Value = 26665
$ws.Range("M25").Value = -26435
$ws.Range("H40").Value = 3457144.8
$ws.Range("I40").Value = 4768980.5
$ws.Range("J40").Value = 13575.25
$ws.Range("K40").Value = 4768980.5
$ws.Range("L40").Value = 13575.25
$ws.Range("M40").Value = -4768844.5
$ws.Range("N40").Value = -13847.25
$ws.Range("H46").Value = 4835273
$ws.Range("I46").Value = 14493287
$ws.Range("J46").Value = 6266.3335
$ws.Range("K46").Value = 14493287
$ws.Range("L46").Value = 6266.3335
$ws.Range("M46").Value = -14493099
$ws.Range("N46").Value = -6642.3335
$ws.Range("H55").Value = 1557.8928
$ws.Range("I55").Value = 934.2381
$ws.Range("J55").Value = 3428.8572
$ws.Range("K55").Value = 934.2381
$ws.Range("L55").Value = 3428.8572
$ws.Range("M55").Value = -761.2381
$ws.Range("N55").Value = -3774.8572
$ws.Range("H61").Value = 12348268
$ws.Range("I61").Value = 18520902
$ws.Range("J61").Value = 3001.6667
$ws.Range("K61").Value = 18520902
$ws.Range("L61").Value = 3001.6667
$ws.Range("M61").Value = -18520700
$ws.Range("N61").Value = -3405.6667
$ws.Range("H82").Value = 75398696
$ws.Range("J82").Value = 2533.3333
$ws.Range("L82").Value = 2533.3333
$ws.Range("N82").Value = -3255.3333
$ws.Range("H85").Value = 75398696
$ws.Range("J85").Value = 2533.3333
$ws.Range("L85").Value = 2533.3333
$ws.Range("N85").Value = -5029.3333
$ws.Range("H113").Value = 12348268
$ws.Range("I113").Value = 18520902
$ws.Range("J113").Value = 3001.6667
$ws.Range("K113").Value = 18520902
$ws.Range("L113").Value = 3001.6667
$ws.Range("M113").Value = -18518732
$ws.Range("N113").Value = -7341.6667
$ws.Range("H122").Value = 7278.2
$ws.Range("I122").Value = 4798.8
$ws.Range("K122").Value = 14396.4
$ws.Range("M122").Value = -11946.4
$ws.Range("H132").Value = 5417.875
$ws.Range("I132").Value = 4752.875
$ws.Range("K132").Value = 14258.625
$ws.Range("M132").Value = -11728.625
$ws.Range("H136").Value = 37990.5
$ws.Range("I136").Value = 53326.25
$ws.Range("J136").Value = 7319
$ws.Range("K136").Value = 159978.75
$ws.Range("L136").Value = 21957
$ws.Range("M136").Value = -157428.75
$ws.Range("N136").Value = -27057

# ---- Sheet WVR: 26 cell updates ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 50002230
$ws.Range("I107").Value = 71431390
$ws.Range("J107").Value = 858
$ws.Range("K107").Value = 214294170
$ws.Range("L107").Value = 2574
$ws.Range("M107").Value = -214292250
$ws.Range("N107").Value = -6414
$ws.Range("H122").Value = 5098.857
$ws.Range("I122").Value = 5282.8335
$ws.Range("K122").Value = 15848.5005
$ws.Range("M122").Value = -13398.5005
$ws.Range("H126").Value = 4584.857
$ws.Range("I126").Value = 4249
$ws.Range("K126").Value = 12747
$ws.Range("M126").Value = -10277
$ws.Range("H132").Value = 45959496
$ws.Range("I132").Value = 90923130
$ws.Range("J132").Value = 995863.75
$ws.Range("K132").Value = 272769390
$ws.Range("L132").Value = 2987591.25
$ws.Range("M132").Value = -272766860
$ws.Range("N132").Value = -2992651.25
$ws.Range("H136").Value = 5381.6777
$ws.Range("I136").Value = 6262.625
$ws.Range("K136").Value = 18787.875
$ws.Range("M136").Value = -16237.875
